$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reordered "IA Control" lists in column A (content/set unchanged, only order) ---
$changes = @(
    [PSCustomObject]@{ Row = 3; Value = 'AU-4,AU-14 (1)' }
    [PSCustomObject]@{ Row = 4; Value = 'CM-6 b,AU-4' }
    [PSCustomObject]@{ Row = 5; Value = 'CM-6 b,SC-5,SC-5 (2)' }
    [PSCustomObject]@{ Row = 6; Value = 'AU-7 a,AC-6 (8),AU-12 (3),AU-7 b,AC-6 (9),AU-8 b,CM-5 (1)' }
    [PSCustomObject]@{ Row = 7; Value = 'AU-7 a,AU-12 (3),CM-6 b,AU-7 b,AU-8 b,AU-12 c,AU-12 a,CM-5 (1)' }
    [PSCustomObject]@{ Row = 11; Value = 'IA-2 (11),IA-2 (12)' }
    [PSCustomObject]@{ Row = 13; Value = 'CM-7 (2),CM-7 (5) (b)' }
    [PSCustomObject]@{ Row = 14; Value = 'CM-7 (2),CM-7 (5) (b)' }
    [PSCustomObject]@{ Row = 16; Value = 'CM-6 b,CM-7 (2)' }
    [PSCustomObject]@{ Row = 21; Value = 'CM-6 b,CM-7 (2)' }
    [PSCustomObject]@{ Row = 22; Value = 'CM-6 b,CM-7 (2)' }
    [PSCustomObject]@{ Row = 37; Value = 'AC-7 a,AC-7 b' }
    [PSCustomObject]@{ Row = 38; Value = 'AC-7 a,AC-7 b' }
    [PSCustomObject]@{ Row = 39; Value = 'AC-7 a,AC-7 b' }
    [PSCustomObject]@{ Row = 40; Value = 'AC-7 a,AC-7 b' }
    [PSCustomObject]@{ Row = 44; Value = 'AU-3 (1),IA-2,IA-8' }
    [PSCustomObject]@{ Row = 45; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 46; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 47; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 48; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 49; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 50; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 51; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 52; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 53; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 54; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 55; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 56; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 57; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 58; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 59; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 60; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 61; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 62; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 63; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 64; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 65; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 66; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 67; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 68; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 69; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 70; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 71; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 72; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 73; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 74; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 75; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 76; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 77; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 78; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 79; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 80; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 81; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 82; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 83; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 84; Value = 'AU-3 (1),MA-4 (1) (a),AU-3' }
    [PSCustomObject]@{ Row = 85; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 86; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 87; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 88; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 89; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 90; Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 91; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 92; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 93; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 94; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 95; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 96; Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3' }
    [PSCustomObject]@{ Row = 97; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 98; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 99; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 100; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 101; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 102; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,AC-2 (4),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 103; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,AC-2 (4),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 104; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,AC-2 (4),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 105; Value = 'AU-3 (1),AU-3,AU-12 c,AC-2 (4),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 106; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,AC-2 (4),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 107; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,AC-2 (4),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 108; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,AC-2 (4),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 109; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,AC-2 (4),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 110; Value = 'AU-3 (1),AU-3,AU-12 c,AU-12 a,AC-2 (4),MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 111; Value = 'AU-14 (1),AU-3 (1),AU-3,AU-12 c,AU-12 a,MA-4 (1) (a)' }
    [PSCustomObject]@{ Row = 121; Value = 'AU-12 c,AU-9' }
    [PSCustomObject]@{ Row = 126; Value = 'AU-12 c,AC-2 (4),CM-5 (1),AC-6 (9)' }
    [PSCustomObject]@{ Row = 128; Value = 'CM-6 b,IA-5 (1) (a),IA-5 (1) (b)' }
    [PSCustomObject]@{ Row = 132; Value = 'SC-8,SC-13,MA-4 c,AC-17 (2)' }
    [PSCustomObject]@{ Row = 133; Value = 'MA-4 (7),SC-10,AC-12,MA-4 e' }
    [PSCustomObject]@{ Row = 134; Value = 'SC-10,AC-12' }
    [PSCustomObject]@{ Row = 135; Value = 'SC-10,AC-12' }
    [PSCustomObject]@{ Row = 136; Value = 'SC-10,AC-11 a' }
    [PSCustomObject]@{ Row = 137; Value = 'AU-7 a,AU-7 (1),AU-14 (1),AU-3 (1),AU-3,AU-6 (4),CM-6 b,AU-12 a,MA-4 (1) (a),CM-5 (1)' }
    [PSCustomObject]@{ Row = 167; Value = 'SC-8,AC-17 (2)' }
    [PSCustomObject]@{ Row = 169; Value = 'SC-8,AC-17 (2)' }
    [PSCustomObject]@{ Row = 171; Value = 'AC-11 b,AC-11 a' }
    [PSCustomObject]@{ Row = 178; Value = 'AU-4 (1),CM-6 b,AU-6 (4)' }
    [PSCustomObject]@{ Row = 179; Value = 'AC-17 (9),AC-17 (1),CM-6 b,CM-7 b' }
    [PSCustomObject]@{ Row = 180; Value = 'CM-6 b,AC-17 (1),CM-7 b' }
    [PSCustomObject]@{ Row = 206; Value = 'SC-28 (1),SC-28' }
    [PSCustomObject]@{ Row = 214; Value = 'AU-12 c,AC-2 (4),AC-6 (9)' }
    [PSCustomObject]@{ Row = 219; Value = 'CM-6 b,IA-2 (5)' }
    [PSCustomObject]@{ Row = 220; Value = 'IA-2 (2),IA-2 (4),IA-2,IA-2 (3),IA-2 (5)' }
    [PSCustomObject]@{ Row = 221; Value = 'IA-2 (2),IA-2 (4),IA-2,IA-2 (3),IA-2 (5)' }
    [PSCustomObject]@{ Row = 222; Value = 'SC-8,AC-18 (1),SC-8 (1)' }
    [PSCustomObject]@{ Row = 224; Value = 'IA-5 (1) (c),IA-7' }
    [PSCustomObject]@{ Row = 225; Value = 'CM-6 b,IA-7' }
    [PSCustomObject]@{ Row = 226; Value = 'CM-6 b,IA-7' }
    [PSCustomObject]@{ Row = 227; Value = 'CM-6 b,IA-7' }
    [PSCustomObject]@{ Row = 229; Value = 'CM-7 a,IA-7' }
    [PSCustomObject]@{ Row = 230; Value = 'SC-13,MA-4 (6)' }
    [PSCustomObject]@{ Row = 231; Value = 'MA-4 (6),AC-17 (2)' }
    [PSCustomObject]@{ Row = 232; Value = 'SC-13,MA-4 (6)' }
    [PSCustomObject]@{ Row = 243; Value = 'CM-6 b,SI-16,SC-2' }
    [PSCustomObject]@{ Row = 257; Value = 'CM-6 b,IA-3' }
    [PSCustomObject]@{ Row = 258; Value = 'CM-6 b,IA-3' }
    [PSCustomObject]@{ Row = 259; Value = 'CM-6 b,IA-3' }
    [PSCustomObject]@{ Row = 260; Value = 'CM-6 b,IA-3' }
    [PSCustomObject]@{ Row = 262; Value = 'AU-5 (1),AU-5 a' }
    [PSCustomObject]@{ Row = 268; Value = 'CM-6 b,IA-2 (2)' }
    [PSCustomObject]@{ Row = 269; Value = 'IA-2 (3),IA-2 (1),IA-2 (2),IA-2 (4)' }
    [PSCustomObject]@{ Row = 274; Value = 'CM-6 b,SC-4' }
    [PSCustomObject]@{ Row = 275; Value = 'SC-4,SC-2' }
    [PSCustomObject]@{ Row = 276; Value = 'SC-4,SC-2' }
    [PSCustomObject]@{ Row = 282; Value = 'CM-6 b,CM-5 (3)' }
    [PSCustomObject]@{ Row = 297; Value = 'IA-2 (11),IA-2 (12)' }
    [PSCustomObject]@{ Row = 298; Value = 'IA-2 (11),IA-2 (1)' }
    [PSCustomObject]@{ Row = 299; Value = 'IA-2 (11),IA-2 (12),IA-2 (1)' }
    [PSCustomObject]@{ Row = 309; Value = 'AU-8 (1) (a),AU-8 (1) (b),AU-8 b' }
    [PSCustomObject]@{ Row = 328; Value = 'AU-12 c,CM-5 (1)' }
    [PSCustomObject]@{ Row = 343; Value = 'CM-7 a,CM-7 b' }
    [PSCustomObject]@{ Row = 344; Value = 'CM-7 a,CM-7 b' }
    [PSCustomObject]@{ Row = 345; Value = 'AC-17 (1),CM-7 b' }
    [PSCustomObject]@{ Row = 346; Value = 'CM-7 a,AC-18 (1)' }
    [PSCustomObject]@{ Row = 347; Value = 'CM-7 a,IA-5 (1) (c),CM-6 b' }
    [PSCustomObject]@{ Row = 358; Value = 'AC-11 b,AC-11 (1)' }
    [PSCustomObject]@{ Row = 367; Value = 'CM-7 a,SI-16' }
    [PSCustomObject]@{ Row = 374; Value = 'CM-7 a,CM-6 b' }
    [PSCustomObject]@{ Row = 375; Value = 'CM-7 a,CM-6 b' }
    [PSCustomObject]@{ Row = 376; Value = 'CM-7 a,CM-6 b' }
    [PSCustomObject]@{ Row = 385; Value = 'CM-6 b,AC-17 (2)' }
    [PSCustomObject]@{ Row = 389; Value = 'SI-6 a,SC-3' }
    [PSCustomObject]@{ Row = 398; Value = 'CM-6 b,SI-16' }
    [PSCustomObject]@{ Row = 401; Value = 'CM-6 b,SC-3' }
    [PSCustomObject]@{ Row = 402; Value = 'CM-6 b,SC-3' }
    [PSCustomObject]@{ Row = 403; Value = 'CM-6 b,SC-3' }
    [PSCustomObject]@{ Row = 450; Value = 'CM-6 b,CM-5 (1)' }
    [PSCustomObject]@{ Row = 451; Value = 'CM-6 b,CM-5 (1)' }
    [PSCustomObject]@{ Row = 524; Value = 'CM-6 b,SC-2' }
    [PSCustomObject]@{ Row = 525; Value = 'CM-6 b,SC-2' }
    [PSCustomObject]@{ Row = 541; Value = 'CM-6 b,SI-2 (2)' }
    [PSCustomObject]@{ Row = 550; Value = 'CM-6 b,SI-2 (2)' }
)

foreach ($item in $changes) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Value
}

# --- Row 208: K208 wording tweak + new M208 "Fix" text ---
$k208 = "To determine if  !authenticate  has not been configured for sudo, run the following command:" + "`n" + " `$ sudo grep -r \!authenticate /etc/sudoers /etc/sudoers.d/ " + "`n" + "The command should return no output." + "`n" + "`n" + "If !authenticate is specified in the sudo config files then this is a finding."
$ws.Cells.Item(208, 11).Value = $k208

$m208 = 'Check that Red Hat Enterprise Linux 9 is not configured to allow users to execute privileged actions without authenticating.' + "`n" + 'Remove any occurrence of "!authenticate" found in "/etc/sudoers" file or files in the "/etc/sudoers.d" directory.' + "`n" + '$ sed -i ''/\!authenticate/ s/^/# /g'' /etc/sudoers /etc/sudoers.d/*'
$ws.Cells.Item(208, 13).Value = $m208

# --- Row 209: K209 wording tweak + new M209 "Fix" text ---
$k209 = "To determine if  NOPASSWD  has been configured for sudo, run the following command:" + "`n" + " `$ sudo grep -ri nopasswd /etc/sudoers /etc/sudoers.d/ " + "`n" + "The command should return no output." + "`n" + "`n" + "If nopasswd is specified in the sudo config files then this is a finding."
$ws.Cells.Item(209, 11).Value = $k209

$m209 = 'Check that Red Hat Enterprise Linux 9 is not configured to allow users to execute privileged actions without authenticating.' + "`n" + 'Remove any occurrence of "NOPASSWD" found in "/etc/sudoers" file or files in the "/etc/sudoers.d" directory.' + "`n" + '$ sed -i ''/NOPASSWD/ s/^/# /g'' /etc/sudoers /etc/sudoers.d/*'
$ws.Cells.Item(209, 13).Value = $m209
